$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1347866666666667
$ws.Range("H2").Value = 0.40436
$ws.Range("I2").Value = 0.03419045085634245
$ws.Range("J2").Value = 0.03419045085634244
$ws.Range("M2").Value = 154.017718
$ws.Range("N2").Value = 462.0531539999999
$ws.Range("O2").Value = 0.9743120958630357
$ws.Range("P2").Value = 0.9743120958630357
$ws.Range("Q2").Value = 20.75953481682666
$ws.Range("R2").Value = 186.83581335144
$ws.Range("S2").Value = 0.03331216983234513
$ws.Range("T2").Value = 0.03331216983234513
$ws.Range("G3").Value = 0.1347866666666667
$ws.Range("H3").Value = 0.40436
$ws.Range("I3").Value = 0.03419045085634245
$ws.Range("J3").Value = 0.03419045085634244
$ws.Range("O3").Value = 0.0142585284421194
$ws.Range("P3").Value = 0.0142585284421194
$ws.Range("Q3").Value = 0.303804518991111
$ws.Range("R3").Value = 2.73424067092
$ws.Range("S3").Value = 0.0004875055159840442
$ws.Range("T3").Value = 0.0004875055159840441
$ws.Range("G4").Value = 0.1347866666666667
$ws.Range("H4").Value = 0.40436
$ws.Range("I4").Value = 0.03419045085634245
$ws.Range("J4").Value = 0.03419045085634244
$ws.Range("O4").Value = 0.011429375694845
$ws.Range("P4").Value = 0.011429375694845
$ws.Range("Q4").Value = 0.2435241476311111
$ws.Range("R4").Value = 2.19171732868
$ws.Range("S4").Value = 0.0003907755080132727
$ws.Range("T4").Value = 0.0003907755080132726
$ws.Range("I5").Value = 0.3318597741685039
$ws.Range("J5").Value = 0.3318597741685039
$ws.Range("M5").Value = 154.017718
$ws.Range("N5").Value = 462.0531539999999
$ws.Range("O5").Value = 0.9743120958630357
$ws.Range("P5").Value = 0.9743120958630357
$ws.Range("Q5").Value = 201.496451892424
$ws.Range("R5").Value = 1813.468067031816
$ws.Range("S5").Value = 0.3233349921027488
$ws.Range("T5").Value = 0.3233349921027487
$ws.Range("I6").Value = 0.3318597741685039
$ws.Range("J6").Value = 0.3318597741685039
$ws.Range("O6").Value = 0.0142585284421194
$ws.Range("P6").Value = 0.0142585284421194
$ws.Range("S6").Value = 0.004731832028776933
$ws.Range("T6").Value = 0.004731832028776932
$ws.Range("I7").Value = 0.3318597741685039
$ws.Range("J7").Value = 0.3318597741685039
$ws.Range("O7").Value = 0.011429375694845
$ws.Range("P7").Value = 0.011429375694845
$ws.Range("S7").Value = 0.003792950036978249
$ws.Range("T7").Value = 0.003792950036978248
$ws.Range("I8").Value = 0.6339497749751537
$ws.Range("J8").Value = 0.6339497749751537
$ws.Range("M8").Value = 154.017718
$ws.Range("N8").Value = 462.0531539999999
$ws.Range("O8").Value = 0.9743120958630357
$ws.Range("P8").Value = 0.9743120958630357
$ws.Range("Q8").Value = 384.9174870788466
$ws.Range("R8").Value = 3464.25738370962
$ws.Range("S8").Value = 0.6176649339279419
$ws.Range("T8").Value = 0.6176649339279419
$ws.Range("I9").Value = 0.6339497749751537
$ws.Range("J9").Value = 0.6339497749751537
$ws.Range("O9").Value = 0.0142585284421194
$ws.Range("P9").Value = 0.0142585284421194
$ws.Range("S9").Value = 0.009039190897358419
$ws.Range("T9").Value = 0.009039190897358419
$ws.Range("I10").Value = 0.6339497749751537
$ws.Range("J10").Value = 0.6339497749751537
$ws.Range("O10").Value = 0.011429375694845
$ws.Range("P10").Value = 0.011429375694845
$ws.Range("R10").Value = 40.63820957388999
$ws.Range("S10").Value = 0.007245650149853477
$ws.Range("T10").Value = 0.007245650149853476
